$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections (2012 price-range counts, column B) ---
# The "35-50" price bucket histogram bars were corrected/shifted.
$ws.Range("B8").Value  = 75
$ws.Range("B9").Value  = 95
$ws.Range("B10").Value = 90
$ws.Range("B11").Value = 75
$ws.Range("B12").Value = 60

# --- Selection now spans the whole data table instead of the stray C33 cell ---
$ws.Range("A1:C32").Select()
